$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": insert a new client row ("EQUISAB S.A.")
# alphabetically before "FARIAS CAICEDO GABRIELA PATRICIA" (row 6),
# pushing every row below it down by one. The sheet's closing summary
# row ("X de 12" -> "X de 13") moves from row 14 to row 15 and its
# counts must be updated to reflect the new total of 13 data rows.
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows.Item(6).Insert()

$ws1.Range("A6").Value = "OFICINA-CATAECSA"
$ws1.Range("B6").Value = "EQUISAB S.A."
$ws1.Range("C6:R6").Value = 0

$ws1.Range("C15").Value = "0 de 13"
$ws1.Range("D15").Value = "0 de 13"
$ws1.Range("E15").Value = "1 de 13"
$ws1.Range("F15").Value = "0 de 13"
$ws1.Range("G15").Value = "0 de 13"
$ws1.Range("H15").Value = "0 de 13"
$ws1.Range("I15").Value = "1 de 13"
$ws1.Range("J15").Value = "0 de 13"
$ws1.Range("K15").Value = "0 de 13"
$ws1.Range("L15").Value = "2 de 13"
$ws1.Range("M15").Value = "1 de 13"
$ws1.Range("N15").Value = "0 de 13"
$ws1.Range("O15").Value = "1 de 13"
$ws1.Range("P15").Value = "0 de 13"
$ws1.Range("Q15").Value = "0 de 13"
$ws1.Range("R15").Value = "0 de 13"

# ------------------------------------------------------------------
# Sheet "VENTA MENSUAL": same new client row, same position. The
# trailing totals row (plain sums) shifts from row 14 to row 15 but
# keeps its values unchanged (the inserted row only adds zeros).
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(6).Insert()

$ws2.Range("A6").Value = "OFICINA-CATAECSA"
$ws2.Range("B6").Value = "EQUISAB S.A."
$ws2.Range("C6:G6").Value = 0
